# Refresh cryptos list with latest price/volume snapshot (GitHub Actions scrape).
# Column D holds numeric-looking text (e.g. "1.00", "0.0000143", "57.131.38") that must
# stay as literal text (matching the original inlineStr cells) rather than being
# reinterpreted by Excel as numbers/dates, so those assignments use a leading apostrophe
# to force text entry and then reset the cell style back to Normal (removing the
# quote-prefix marker) so no visible formatting changes are introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$sub3 = [char]0x2083  # subscript-3 glyph used in small-number notation (e.g. 0.0₃0867)

$ws.Range("D2").Value = "'57.131.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -5.20%  "
$ws.Range("D3").Value = "'3.062.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -7.05%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'511.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.59%  "
$ws.Range("D6").Value = "'127.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.30%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'3.057.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.30%  "
$ws.Range("D9").Value = "'0.430"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.03%  "
$ws.Range("D10").Value = "'7.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.90%  "
$ws.Range("D11").Value = "'0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.57%  "
$ws.Range("D12").Value = "'0.358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.62%  "
$ws.Range("E13").Value = "  -6.47%  "
$ws.Range("D14").Value = "'0.127"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "'24.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.68%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'3.089.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.49%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'54.660.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.24%  "
$ws.Range("D18").Value = "'0.0000143"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -12.00%  "
$ws.Range("D19").Value = "'5.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.29%  "
$ws.Range("D20").Value = "'12.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.61%  "
$ws.Range("D21").Value = "'7.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.77%  "
$ws.Range("D22").Value = "'326.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -12.22%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'65.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.22%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.486"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.11%  "
$ws.Range("D26").Value = "'0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.43%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'0.0" + $sub3 + "0867"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.46%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'6.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.50%  "
$ws.Range("D31").Value = "'1.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").Value = "'1.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.05%  "
$ws.Range("D33").Value = "'6.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.59%  "
$ws.Range("D34").Value = "'20.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.34%  "
$ws.Range("D35").Value = "'156.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.73%  "
$ws.Range("D36").Value = "'4.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.66%  "
$ws.Range("D37").Value = "'5.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.88%  "
$ws.Range("D38").Value = "'1.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -13.42%  "
$ws.Range("D39").Value = "'3.112.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.56%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0654"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.42%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'22.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.60%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'36.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -12.68%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'0.663"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.82%  "
$ws.Range("D45").Value = "'1.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.23%  "
$ws.Range("D46").Value = "'3.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.02%  "
$ws.Range("D47").Value = "'2.216.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("D48").Value = "'1.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.49%  "
$ws.Range("D49").Value = "'5.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.95%  "
$ws.Range("D50").Value = "'19.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.62%  "
$ws.Range("D51").Value = "'0.0225"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.71%  "
